# Weekly update: insert a new price record for Puerro (Vega Modelo de Temuco)
# as the new first data row (row 118), pushing all existing records down by
# one row (old row 118 -> 119, ..., old row 229 -> 230).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 118; Excel shifts rows 118:229 down to 119:230.
$ws.Rows(118).Insert()

# Populate the new row 118 with the new weekly price entry.
$ws.Range("A118").Value = 10
$ws.Range("B118").Value = "Vega Modelo de Temuco"
$ws.Range("C118").Value = "La Araucanía"
$ws.Range("D118").Value = 44790
$ws.Range("E118").Value = 9
$ws.Range("F118").Value = 100112005
$ws.Range("G118").Value = "Puerro"
$ws.Range("H118").Value = "Azul de Maquehue"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 65
$ws.Range("K118").Value = 16000
$ws.Range("L118").Value = 16000
$ws.Range("M118").Value = 16000
$ws.Range("N118").Value = "`$/docena de paquetes"
$ws.Range("O118").Value = "Provincia de Cautín"
$ws.Range("P118").Value = 1333
$ws.Range("Q118").Value = 12
$ws.Range("R118").Value = "Hortaliza"
